# "Game loop finished (finally)" / "No audio at all"
#
# This workbook tracks items found during a game. Every remaining item
# has now been located, so the "Done?" column is set to "Y" for every
# row that didn't already have it (rows 72-76: Waifu, Wallet, Wife,
# Wine, Wristwatch). Also corrects a mis-typed item name ("Plot" with a
# stray curly closing quote) to '"Plot"' (opening straight quote +
# closing curly quote).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix item text in B50 (id 49): Plot" -> "Plot"
$ws.Range("B50").Value = '"Plot' + [char]0x201D

# Mark the last five items (rows 72-76) as done
$ws.Range("E72:E76").Value = "Y"

# Restore the scroll/selection state left by the editing session
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D41").Select()
